# Weekly fruit/vegetable price update: shift existing rows 12-21 down to
# 13-22 (one week older each) and insert a new latest-week row at row 12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 is brand new -- seed the columns that stay constant across all
# rows of this series (copied from the template row).
$ws.Cells.Item(22, 1).Value = 10
$ws.Cells.Item(22, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(22, 3).Value = 'La Araucanía'
$ws.Cells.Item(22, 5).Value = 9
$ws.Cells.Item(22, 6).Value = 300000001
$ws.Cells.Item(22, 7).Value = 'Rabanito'
$ws.Cells.Item(22, 8).Value = 'Sin especificar'
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(22, 17).Value = 12
$ws.Cells.Item(22, 18).Value = 'Hortaliza'

# Row 12: Fecha=44413, Volumen=40, Origen=Provincia de Cautín
$ws.Cells.Item(12, 4).Value = 44413
$ws.Cells.Item(12, 10).Value = 40
$ws.Cells.Item(12, 11).Value = 7000
$ws.Cells.Item(12, 12).Value = 7000
$ws.Cells.Item(12, 13).Value = 7000
$ws.Cells.Item(12, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(12, 16).Value = 583

# Row 13: Fecha=44410, Volumen=40, Origen=Provincia de Cautín
$ws.Cells.Item(13, 4).Value = 44410
$ws.Cells.Item(13, 10).Value = 40
$ws.Cells.Item(13, 11).Value = 7000
$ws.Cells.Item(13, 12).Value = 7000
$ws.Cells.Item(13, 13).Value = 7000
$ws.Cells.Item(13, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(13, 16).Value = 583

# Row 14: Fecha=44327, Volumen=30, Origen=Provincia de Cautín
$ws.Cells.Item(14, 4).Value = 44327
$ws.Cells.Item(14, 10).Value = 30
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 6000
$ws.Cells.Item(14, 13).Value = 6000
$ws.Cells.Item(14, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(14, 16).Value = 500

# Row 15: Fecha=44196, Volumen=20, Origen=Provincia de Cautín
$ws.Cells.Item(15, 4).Value = 44196
$ws.Cells.Item(15, 10).Value = 20
$ws.Cells.Item(15, 11).Value = 5000
$ws.Cells.Item(15, 12).Value = 5000
$ws.Cells.Item(15, 13).Value = 5000
$ws.Cells.Item(15, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(15, 16).Value = 417

# Row 16: Fecha=44369, Volumen=20, Origen=Región Metropolitana
$ws.Cells.Item(16, 4).Value = 44369
$ws.Cells.Item(16, 10).Value = 20
$ws.Cells.Item(16, 11).Value = 4000
$ws.Cells.Item(16, 12).Value = 4000
$ws.Cells.Item(16, 13).Value = 4000
$ws.Cells.Item(16, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(16, 16).Value = 333

# Row 17: Fecha=44195, Volumen=55, Origen=Provincia de Cautín
$ws.Cells.Item(17, 4).Value = 44195
$ws.Cells.Item(17, 10).Value = 55
$ws.Cells.Item(17, 11).Value = 5000
$ws.Cells.Item(17, 12).Value = 5000
$ws.Cells.Item(17, 13).Value = 5000
$ws.Cells.Item(17, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(17, 16).Value = 417

# Row 18: Fecha=44186, Volumen=50, Origen=Provincia de Cautín
$ws.Cells.Item(18, 4).Value = 44186
$ws.Cells.Item(18, 10).Value = 50
$ws.Cells.Item(18, 11).Value = 5000
$ws.Cells.Item(18, 12).Value = 5000
$ws.Cells.Item(18, 13).Value = 5000
$ws.Cells.Item(18, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(18, 16).Value = 417

# Row 19: Fecha=44211, Volumen=65, Origen=Provincia de Cautín
$ws.Cells.Item(19, 4).Value = 44211
$ws.Cells.Item(19, 10).Value = 65
$ws.Cells.Item(19, 11).Value = 5000
$ws.Cells.Item(19, 12).Value = 5000
$ws.Cells.Item(19, 13).Value = 5000
$ws.Cells.Item(19, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(19, 16).Value = 417

# Row 20: Fecha=44301, Volumen=50, Origen=Provincia de Cautín
$ws.Cells.Item(20, 4).Value = 44301
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 6000
$ws.Cells.Item(20, 12).Value = 6000
$ws.Cells.Item(20, 13).Value = 6000
$ws.Cells.Item(20, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(20, 16).Value = 500

# Row 21: Fecha=44326, Volumen=50, Origen=Provincia de Cautín
$ws.Cells.Item(21, 4).Value = 44326
$ws.Cells.Item(21, 10).Value = 50
$ws.Cells.Item(21, 11).Value = 6000
$ws.Cells.Item(21, 12).Value = 6000
$ws.Cells.Item(21, 13).Value = 6000
$ws.Cells.Item(21, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(21, 16).Value = 500

# Row 22: Fecha=44179, Volumen=40, Origen=Provincia de Cautín
$ws.Cells.Item(22, 4).Value = 44179
$ws.Cells.Item(22, 10).Value = 40
$ws.Cells.Item(22, 11).Value = 6000
$ws.Cells.Item(22, 12).Value = 6000
$ws.Cells.Item(22, 13).Value = 6000
$ws.Cells.Item(22, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(22, 16).Value = 500

# Column D carries a date/time number format (numFmtId 165) on every data
# row; make sure the newly written D12:D22 values keep using it.
$ws.Range('D12:D22').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
